$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text edits (masthead volume/number + date range) ---
$ws.Range("A8").Value = "Volume 32   Number  12"
$ws.Range("C9").Value = "Report Covering the Week  3/17/2025  Through  3/23/2025"

# --- Cell type changes (number <-> "N/A"-style shared text) ---
# Done via Copy() from a stable, diff-unaffected template cell (row 31) to pick up
# the correct style index, then overwriting the value.
$ws.Range("C31").Copy($ws.Range("F14"))
$ws.Range("C31").Copy($ws.Range("C22"))
$ws.Range("G31").Copy($ws.Range("D29"))
$ws.Range("D29").Value = 1
$ws.Range("H31").Copy($ws.Range("E29"))
$ws.Range("E29").Value = -100
$ws.Range("C31").Copy($ws.Range("F29"))
$ws.Range("G31").Copy($ws.Range("D30"))
$ws.Range("D30").Value = 1
$ws.Range("H31").Copy($ws.Range("E30"))
$ws.Range("E30").Value = -100
$ws.Range("C31").Copy($ws.Range("F30"))
$ws.Range("G31").Copy($ws.Range("D33"))
$ws.Range("D33").Value = 1
$ws.Range("H31").Copy($ws.Range("E33"))
$ws.Range("E33").Value = -100
$ws.Range("G31").Copy($ws.Range("G33"))
$ws.Range("G33").Value = 1
$ws.Range("H31").Copy($ws.Range("H33"))
$ws.Range("H33").Value = 0
$ws.Range("G31").Copy($ws.Range("J33"))
$ws.Range("J33").Value = 1
$ws.Range("H31").Copy($ws.Range("K33"))
$ws.Range("K33").Value = 200

# --- Plain numeric overwrites (style/format unchanged) ---
# Row 14
$ws.Range("N14").Value = -95.833333333333
# Row 15
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 200
$ws.Range("F15").Value = 7
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 15
$ws.Range("J15").Value = 14
$ws.Range("K15").Value = 7.142857142857
$ws.Range("L15").Value = -11.764705882352
$ws.Range("N15").Value = -25
# Row 16
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = -10
$ws.Range("G16").Value = 51
$ws.Range("H16").Value = -29.411764705882
$ws.Range("I16").Value = 83
$ws.Range("J16").Value = 140
$ws.Range("K16").Value = -40.714285714285
$ws.Range("L16").Value = -40.714285714285
$ws.Range("M16").Value = -37.593984962406
$ws.Range("N16").Value = -86.437908496732
# Row 17
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = -6.25
$ws.Range("F17").Value = 68
$ws.Range("G17").Value = 78
$ws.Range("H17").Value = -12.820512820512
$ws.Range("I17").Value = 208
$ws.Range("J17").Value = 237
$ws.Range("K17").Value = -12.236286919831
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 33.333333333333
$ws.Range("N17").Value = -29.96632996633
# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -50
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = -5.263157894736
$ws.Range("I18").Value = 48
$ws.Range("J18").Value = 55
$ws.Range("K18").Value = -12.727272727272
$ws.Range("L18").Value = -52.475247524752
$ws.Range("M18").Value = -47.252747252747
$ws.Range("N18").Value = -87.061994609164
# Row 19
$ws.Range("C19").Value = 22
$ws.Range("D19").Value = 21
$ws.Range("E19").Value = 4.761904761904
$ws.Range("F19").Value = 67
$ws.Range("G19").Value = 74
$ws.Range("H19").Value = -9.459459459459
$ws.Range("I19").Value = 175
$ws.Range("J19").Value = 210
$ws.Range("K19").Value = -16.666666666666
$ws.Range("L19").Value = -28.27868852459
$ws.Range("M19").Value = 17.44966442953
$ws.Range("N19").Value = 5.421686746987
# Row 20
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 35
$ws.Range("H20").Value = -5.405405405405
$ws.Range("I20").Value = 83
$ws.Range("J20").Value = 119
$ws.Range("K20").Value = -30.252100840336
$ws.Range("L20").Value = -35.15625
$ws.Range("M20").Value = 59.615384615384
$ws.Range("N20").Value = -87.918486171761
# Row 21
$ws.Range("C21").Value = 58
$ws.Range("D21").Value = 63
$ws.Range("E21").Value = -7.936507936507
$ws.Range("F21").Value = 231
$ws.Range("G21").Value = 266
$ws.Range("H21").Value = -13.157894736842
$ws.Range("I21").Value = 613
$ws.Range("J21").Value = 778
$ws.Range("K21").Value = -21.208226221079
$ws.Range("L21").Value = -27.110582639714
$ws.Range("M21").Value = 1.827242524916
$ws.Range("N21").Value = -71.841984382177
# Row 22
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 200
$ws.Range("L22").Value = -18.181818181818
# Row 23
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 0
$ws.Range("G23").Value = 17
$ws.Range("H23").Value = 47.058823529411
$ws.Range("I23").Value = 59
$ws.Range("J23").Value = 54
$ws.Range("K23").Value = 9.259259259259
$ws.Range("L23").Value = -33.707865168539
$ws.Range("M23").Value = 78.787878787878
# Row 24
$ws.Range("C24").Value = 41
$ws.Range("D24").Value = 52
$ws.Range("E24").Value = -21.153846153846
$ws.Range("F24").Value = 187
$ws.Range("G24").Value = 183
$ws.Range("H24").Value = 2.185792349726
$ws.Range("I24").Value = 563
$ws.Range("J24").Value = 508
$ws.Range("K24").Value = 10.826771653543
$ws.Range("L24").Value = 12.15139442231
$ws.Range("M24").Value = 67.062314540059
# Row 25
$ws.Range("C25").Value = 19
$ws.Range("D25").Value = 27
$ws.Range("E25").Value = -29.629629629629
$ws.Range("F25").Value = 66
$ws.Range("G25").Value = 80
$ws.Range("H25").Value = -17.5
$ws.Range("I25").Value = 205
$ws.Range("J25").Value = 201
$ws.Range("K25").Value = 1.990049751243
$ws.Range("L25").Value = 14.525139664804
# Row 26
$ws.Range("C26").Value = 20
$ws.Range("D26").Value = 15
$ws.Range("E26").Value = 33.333333333333
$ws.Range("F26").Value = 115
$ws.Range("G26").Value = 106
$ws.Range("H26").Value = 8.490566037735
$ws.Range("I26").Value = 302
$ws.Range("J26").Value = 322
$ws.Range("K26").Value = -6.211180124223
$ws.Range("L26").Value = 26.890756302521
$ws.Range("M26").Value = -18.817204301075
# Row 27
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = -20
$ws.Range("I27").Value = 21
$ws.Range("J27").Value = 25
$ws.Range("K27").Value = -16
$ws.Range("L27").Value = 10.526315789473
# Row 28
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 50
$ws.Range("F28").Value = 13
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = 44.444444444444
$ws.Range("I28").Value = 25
$ws.Range("J28").Value = 21
$ws.Range("K28").Value = 19.047619047619
$ws.Range("L28").Value = 25
# Row 29
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = -100
$ws.Range("J29").Value = 13
$ws.Range("K29").Value = -23.076923076923
$ws.Range("L29").Value = -41.176470588235
$ws.Range("M29").Value = -33.333333333333
$ws.Range("N29").Value = -87.951807228915
# Row 30
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = -100
$ws.Range("J30").Value = 12
$ws.Range("K30").Value = -33.333333333333
$ws.Range("L30").Value = -42.857142857142
$ws.Range("M30").Value = -38.461538461538
$ws.Range("N30").Value = -89.873417721519
# Row 33
$ws.Range("F33").Value = 1
$ws.Range("I33").Value = 3

# --- Column width (E column bestFit widened) ---
$ws.Columns("E").ColumnWidth = $ws.Columns("H").ColumnWidth
